$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously held two sample applicants (rows 2 and 3) - remove that
# sample data, keeping only the header row behind.
$ws.Rows("2:3").Delete()

# Add a new "AvailableOnDayX" field, inserted right before "GotSomethingElseToSay"
# (this was column K before the new leading Id column is added below).
$ws.Columns("K").Insert()
$ws.Range("K1").Value = "AvailableOnDayX"

# Add a new leading "Id" column to hold an auto incrementing applicant id.
$ws.Columns("A").Insert()
$ws.Range("A1").Value = "Id"

# Re-apply/tidy up the column widths for the new, wider form layout.
$ws.Columns("A").ColumnWidth = 16.333333333333332
$ws.Columns("B").ColumnWidth = 22.833333333333332
$ws.Columns("C").ColumnWidth = 19.833333333333332
$ws.Columns("D").ColumnWidth = 20.666666666666668
$ws.Columns("E").ColumnWidth = 30.666666666666668
$ws.Columns("F").ColumnWidth = 9.666666666666666
$ws.Columns("G").ColumnWidth = 25.5
$ws.Columns("H").ColumnWidth = 16.666666666666668
$ws.Columns("I").ColumnWidth = 13.5
$ws.Columns("J").ColumnWidth = 9.666666666666666
$ws.Columns("K").ColumnWidth = 11.833333333333334
$ws.Columns("L").ColumnWidth = 16.333333333333332
$ws.Columns("M").ColumnWidth = 32.5
$ws.Columns("N").ColumnWidth = 35.5
$ws.Columns("P").ColumnWidth = 13.333333333333334

$ws.Range("P1").Select()

$wb.Save()
